$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VTRS")

# Row 4 (Inventory)
$ws.Range("B4").Value = 5472000000.0
$ws.Range("C4").Value = 3022000000.0
$ws.Range("D4").Value = 2786000000.0
$ws.Range("E4").Value = 2640000000.0
$ws.Range("F4").Value = 2671000000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 1346000000.0
$ws.Range("C14").Value = 901000000.0
$ws.Range("D14").Value = 801000000.0
$ws.Range("E14").Value = 875000000.0
$ws.Range("F14").Value = 1062000000.0

# Row 21 (Long Term Tax Liability (Deferred))
$ws.Range("B21").Value = 976000000.0
$ws.Range("C21").Value = 759000000.0
$ws.Range("D21").Value = 793000000.0
$ws.Range("E21").Value = 837000000.0
$ws.Range("F21").Value = 924000000.0
